# The deck's design theme was changed from the "Integral" (Red Violet)
# theme to the built-in "Office Theme" palette. This replicates that
# action by rewriting the 12 theme colors (and, defensively, the
# theme/color-scheme display names) on the presentation's theme so the
# underlying ppt/theme/*.xml color scheme matches the standard Office
# theme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

function ToVbaRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Index -> (scheme slot, target Office Theme RGB)
$officeColors = @(
    @(1,  0x00, 0x00, 0x00),   # dk1
    @(2,  0xFF, 0xFF, 0xFF),   # lt1
    @(3,  0x44, 0x54, 0x6A),   # dk2
    @(4,  0xE7, 0xE6, 0xE6),   # lt2
    @(5,  0x5B, 0x9B, 0xD5),   # accent1
    @(6,  0xED, 0x7D, 0x31),   # accent2
    @(7,  0xA5, 0xA5, 0xA5),   # accent3
    @(8,  0xFF, 0xC0, 0x00),   # accent4
    @(9,  0x44, 0x72, 0xC4),   # accent5
    @(10, 0x70, 0xAD, 0x47),   # accent6
    @(11, 0x05, 0x63, 0xC1),   # hlink
    @(12, 0x95, 0x4F, 0x72)    # folHlink
)

foreach ($entry in $officeColors) {
    $idx = $entry[0]
    $r = $entry[1]
    $g = $entry[2]
    $b = $entry[3]
    $colorScheme.Item($idx).RGB = ToVbaRGB $r $g $b
}

# Best-effort: rename the theme / color scheme to match the standard
# "Office Theme" / "Office" naming (no-op on hosts that keep these
# read-only).
try { $theme.Name = "Office Theme" } catch {}
try { $colorScheme.Name = "Office" } catch {}
